$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New forecast-year header (column Z = "2025/26") ---
$ws.Cells.Item(1, 26).Value = "2025/26"

# --- Fix number format on the "Gross tax revenue as a percentage of GDP" row (2022 block) ---
$ws.Range("V266:Y266").NumberFormat = "_(* #,##0.000_);_(* \(#,##0.000\);_(* ""-""??_);_(@_)"

# --- Append the 2023 forecast block (rows 269-283), mirroring the 2022 block (rows 254-268) ---

function Set-Row($r, $catText, $catWrap, $dVal, $dWrap, $w, $x, $y, $z, $fmt) {
    $ws.Cells.Item($r, 1).Value = "Budget"
    $ws.Cells.Item($r, 2).Value = 2023
    $ws.Cells.Item($r, 3).Value = $catText
    $ws.Cells.Item($r, 4).Value = $dVal
    if ($catWrap) {
        $ws.Cells.Item($r, 3).WrapText = $true
    }
    if ($dWrap) {
        $ws.Cells.Item($r, 4).WrapText = $true
    }
    $ws.Cells.Item($r, 23).Value = $w
    $ws.Cells.Item($r, 24).Value = $x
    $ws.Cells.Item($r, 25).Value = $y
    $ws.Cells.Item($r, 26).Value = $z
    $ws.Range($ws.Cells.Item($r, 23), $ws.Cells.Item($r, 26)).NumberFormat = $fmt
}

$fmtInt = "#,##0"
$fmtDec2 = "#,##0.00"
$fmtDec3 = "_(* #,##0.000_);_(* \(#,##0.000\);_(* ""-""??_);_(@_)"
$fmtBuoy = "0.00"

Set-Row 269 "Taxes on income and profits" $false 1 $false 989877 1021213 1089123 1172033 $fmtInt
Set-Row 270 "Personal income tax" $false 2 $false 601649 640300 696624 752627 $fmtInt
Set-Row 271 "Corporate income tax" $false 3 $false 344944 336119 345434 369477 $fmtInt
Set-Row 272 "Taxes on payroll and workforce" $false 4 $false 21238 23027 24816 26846 $fmtInt
Set-Row 273 "Taxes on property" $false 5 $false 22656 23863 25384 27040 $fmtInt
Set-Row 274 "Domestic taxes on goods and services" $true 6 $true 581871 642765 687208 731032 $fmtInt
Set-Row 275 "Value-added tax" $false 7 $true 426283 471477 505409 537868 $fmtInt
Set-Row 276 "Taxes on international trade and transactions" $true 8 $true 76535 76588 81195 86506 $fmtInt
Set-Row 277 "Gross tax revenue" $false 9 $false 1692177 1787456 1907727 2043456 $fmtInt
Set-Row 278 "Departmental revenue" $false 10 $false 55078 51583 46859 44310 $fmtInt
Set-Row 279 "Less: SACU payments" $false 11 $false -43683 -79811 -86505 -80059 $fmtInt
Set-Row 280 "Main budget revenue" $false 12 $false 1703571 1759229 1868080 2007707 $fmtInt
Set-Row 281 "Gross tax revenue as a percentage of GDP" $true 13 $true 0.254 0.255 0.256 0.257 $fmtDec3
Set-Row 282 "GDP (R billion)" $false 14 $true 6651.3 7005.7 7452.4 7938.5 $fmtDec2
Set-Row 283 "Tax buoyancy" $false 15 $true 1.42 1.06 1.06 1.09 $fmtBuoy

# --- New column Z width ---
$ws.Columns.Item(26).ColumnWidth = 11.43

# --- Frozen-pane / selection state matches a scroll back to the top-left of the data ---
$ws.Range("A2").Select()
